$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.216.90'
$ws.Range("D3").Value = '3.914.33'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '487.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.72%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.731'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.48%  '
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000345'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.24%  '
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.85'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.83%  '
$ws.Range("D14").Value = '4.536.60'
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '3.939.74'
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.136'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = '  -2.30%  '
$ws.Range("D20").Value = '68.322.28'
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '441.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.82%  '
$ws.Range("E23").Value = '  +2.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +17.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +12.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '725.26'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("E32").Value = '  -0.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.91'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.42%  '
$ws.Range("E34").Value = '  +17.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '42.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '61.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.44%  '
$ws.Range("D37").Value = '0.0₃0865'
$ws.Range("E37").Value = '  +7.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.411'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +22.54%  '
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +16.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0481'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("E44").Value = '  +4.23%  '
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0353'
$ws.Range("E49").Value = '  +33.43%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '146.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.52%  '
